# Apply the commit: add building-function values 1431, 1432 and 1650
# to the classification table on sheet "AX_Gebaeudefunktion"
# (file AX_Bauwerksfunktion_BauwerkOderAnlageFuerSportFreizeitUndErholung.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$urlValue = "http://inspire.ec.europa.eu/codelist/BuildingNatureValue/stadium"
$sonstigesValue = "{{project:BUILDINGNATURE}}sonstiges"

# Row 3 (new): Zuschauertribuene, ueberdacht / 1431
$ws.Range("A3").Value = "Zuschauertribüne, überdacht"
$ws.Range("B3").Value = 1431
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = $urlValue

# Row 4 (new): Zuschauertribuene, nicht ueberdacht / 1432
$ws.Range("A4").Value = "Zuschauertribüne, nicht überdacht"
$ws.Range("B4").Value = 1432
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = $urlValue

# Row 5: Stadion / 1440 (shifted down from former row 3)
$ws.Range("A5").Value = "Stadion"
$ws.Range("B5").Value = 1440
$ws.Range("C5").Value = $urlValue

# Row 6: Stadion, ueberdacht / 1441 (shifted down from former row 4)
$ws.Range("A6").Value = "Stadion, überdacht"
$ws.Range("B6").Value = 1441
$ws.Range("C6").Value = $urlValue
$ws.Range("C6").Style = "Hyperlink"

# Row 7: Stadion, nicht ueberdacht / 1442 (shifted down from former row 5)
$ws.Range("A7").Value = "Stadion, nicht überdacht"
$ws.Range("B7").Value = 1442
$ws.Range("C7").Value = $urlValue
$ws.Range("C7").Style = "Hyperlink"

# Row 8: Sprungschanze (Anlauf) / 1470 (shifted down from former row 6)
$ws.Range("A8").Value = "Sprungschanze (Anlauf)"
$ws.Range("B8").Value = 1470
$ws.Range("C8").Value = $urlValue
$ws.Range("C8").Style = "Hyperlink"

# Row 9: Gradierwerk / 1490 (shifted down from former row 7)
$ws.Range("A9").Value = "Gradierwerk"
$ws.Range("B9").Value = 1490
$ws.Range("C9").Value = $sonstigesValue

# Row 10 (new): Wassersportanlage / 1650
$ws.Range("A10").Value = "Wassersportanlage"
$ws.Range("B10").Value = 1650
$ws.Range("C10").Value = $sonstigesValue
# match the plain (non-hyperlink) wrapped-text look used by the "sonstiges" cell above (C9)
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wb.Save()
